# Add Input [Sample Name] / Output [Sample Name] columns to the
# SugarExtraction annotation table (RPTU PP Sugar Extraction template).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SugarExtraction")
$lo = $ws.ListObjects.Item(1)

# Insert a blank column before column A; this shifts the table and all
# of its data one column to the right (B:I) without re-typing any
# cells, so the original stored value types are preserved exactly
# (e.g. the numeric-looking "3"/"95" text stays text).
$ws.Columns.Item(1).Insert()
$lo.Resize($ws.Range("B1:I2"))

# Grow the table so it also covers the new leading column (A) and a
# new trailing column (J).
$lo.Resize($ws.Range("A1:J2"))

# Populate the two new data cells.
$ws.Range("A2").Value = "PlantHarvest"
$ws.Range("J2").Value = "PlantSugarExtract"

# The table's cached column-name list only refreshes for header cells
# written through the table object itself, so rewrite every header (in
# its final left-to-right order) via HeaderRowRange. This keeps
# xl/tables/table1.xml's column names/order in sync with the
# worksheet's actual header cells.
$headers = @( `
    "Input [Sample Name]", `
    "Parameter [Vortex Mixer]", `
    "Unit", `
    "Term Source REF (NCIT:C29544)", `
    "Term Accession Number (NCIT:C29544)", `
    "Parameter [Temperature]", `
    "Unit ", `
    "Term Source REF (NCIT:C25206)", `
    "Term Accession Number (NCIT:C25206)", `
    "Output [Sample Name]" `
)
for ($i = 1; $i -le $headers.Count; $i++) {
    $lo.HeaderRowRange.Cells.Item(1, $i).Value = $headers[$i - 1]
}

# Match the author's final UI state: SugarExtraction is the active sheet.
$ws.Activate()
